$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 2.4
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("U2").Value = 2.1
$ws.Range("V2").Value = 1.67
$ws.Range("W2").Value = 6
$ws.Range("X2").Value = 7.5
$ws.Range("AC2").Value = 8
$ws.Range("AE2").Value = 19
$ws.Range("AF2").Value = 67
$ws.Range("AI2").Value = 21
$ws.Range("AJ2").Value = 15
